$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wrap each theme-name value in column B (rows 2-20) with single quotes,
# leaving the header row (B1 = "description_theme") untouched.
for ($r = 2; $r -le 20; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $current = $cell.Value2
    # A single leading apostrophe is Excel's "force text" prefix and is
    # stripped from the stored value, so it must be doubled to end up with
    # one literal leading apostrophe in the saved cell content.
    $cell.Value = "''" + $current + "'"
    # Writing a string that starts with an apostrophe also marks the cell
    # with the "quote prefix" cell style; restore the plain Normal style so
    # the cell keeps its original (unstyled) formatting.
    $cell.Style = "Normal"
}

# Move the active selection from B21 to F5.
$ws.Range("F5").Select()
